$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: copy formatting (styles) from row 195 into the new rows 196-214,
# and set the row height to match the existing data rows.
$ws.Range("A195:M195").Copy()
$ws.Range("A196:M214").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("196:214").RowHeight = 16

# Step 2: a few text columns hold purely-numeric-looking strings (e.g. document
# IDs like "20639"). Mark those specific cells as Text format first so Excel
# keeps them as strings instead of silently converting them to numbers.
$ws.Range("D196,D197,I197,D198,I198,D199,I199,D200,I200,D201,I201,D202,I202,D203,I203,D204,D205,D206,D207,D208,D209,D210,D211,D212,D213,D214").NumberFormat = "@"

# Step 3: write the cell values for each new row.
# Row 196
$ws.Range("A196").Value = "●"
$ws.Range("D196").Value = "20639"
$ws.Range("E196").Value = "B"
$ws.Range("F196").Value = "2: 5200"
$ws.Range("G196").Value = "2: 5208"
$ws.Range("H196").Value = 0
$ws.Range("I196").Value = "Meropenem"
$ws.Range("J196").Value = 9
$ws.Range("K196").Value = 0.066074443873430741
$ws.Range("L196").Value = "chen"
$ws.Range("M196").Value = "1/30/19 17:57:04"

# Row 197
$ws.Range("A197").Value = "●"
$ws.Range("D197").Value = "21149"
$ws.Range("E197").Value = "MIC"
$ws.Range("F197").Value = "2: 6486"
$ws.Range("G197").Value = "2: 6489"
$ws.Range("H197").Value = 0
$ws.Range("I197").Value = ".256"
$ws.Range("J197").Value = 4
$ws.Range("K197").Value = 0.020821404403727033
$ws.Range("L197").Value = "chen"
$ws.Range("M197").Value = "1/30/19 17:58:03"

# Row 198
$ws.Range("A198").Value = "●"
$ws.Range("D198").Value = "21149"
$ws.Range("E198").Value = "MIC"
$ws.Range("F198").Value = "2: 6509"
$ws.Range("G198").Value = "2: 6510"
$ws.Range("H198").Value = 0
$ws.Range("I198").Value = "64"
$ws.Range("J198").Value = 2
$ws.Range("K198").Value = 0.010410702201863516
$ws.Range("L198").Value = "chen"
$ws.Range("M198").Value = "1/30/19 17:58:12"

# Row 199
$ws.Range("A199").Value = "●"
$ws.Range("D199").Value = "21149"
$ws.Range("E199").Value = "MIC"
$ws.Range("F199").Value = "2: 6530"
$ws.Range("G199").Value = "2: 6533"
$ws.Range("H199").Value = 0
$ws.Range("I199").Value = ".256"
$ws.Range("J199").Value = 4
$ws.Range("K199").Value = 0.020821404403727033
$ws.Range("L199").Value = "chen"
$ws.Range("M199").Value = "1/30/19 17:58:19"

# Row 200
$ws.Range("A200").Value = "●"
$ws.Range("D200").Value = "21149"
$ws.Range("E200").Value = "MIC"
$ws.Range("F200").Value = "2: 6584"
$ws.Range("G200").Value = "2: 6587"
$ws.Range("H200").Value = 0
$ws.Range("I200").Value = ".256"
$ws.Range("J200").Value = 4
$ws.Range("K200").Value = 0.020821404403727033
$ws.Range("L200").Value = "chen"
$ws.Range("M200").Value = "1/30/19 17:58:23"

# Row 201
$ws.Range("A201").Value = "●"
$ws.Range("D201").Value = "21149"
$ws.Range("E201").Value = "MIC"
$ws.Range("F201").Value = "2: 6616"
$ws.Range("G201").Value = "2: 6619"
$ws.Range("H201").Value = 0
$ws.Range("I201").Value = ".256"
$ws.Range("J201").Value = 4
$ws.Range("K201").Value = 0.020821404403727033
$ws.Range("L201").Value = "chen"
$ws.Range("M201").Value = "1/30/19 17:58:31"

# Row 202
$ws.Range("A202").Value = "●"
$ws.Range("D202").Value = "21149"
$ws.Range("E202").Value = "MIC"
$ws.Range("F202").Value = "2: 6370"
$ws.Range("G202").Value = "2: 6371"
$ws.Range("H202").Value = 0
$ws.Range("I202").Value = "32"
$ws.Range("J202").Value = 2
$ws.Range("K202").Value = 0.010410702201863516
$ws.Range("L202").Value = "chen"
$ws.Range("M202").Value = "1/30/19 17:58:37"

# Row 203
$ws.Range("A203").Value = "●"
$ws.Range("D203").Value = "21149"
$ws.Range("E203").Value = "MIC"
$ws.Range("F203").Value = "2: 6462"
$ws.Range("G203").Value = "2: 6463"
$ws.Range("H203").Value = 0
$ws.Range("I203").Value = "64"
$ws.Range("J203").Value = 2
$ws.Range("K203").Value = 0.010410702201863516
$ws.Range("L203").Value = "chen"
$ws.Range("M203").Value = "1/30/19 17:59:00"

# Row 204
$ws.Range("A204").Value = "●"
$ws.Range("D204").Value = "23107"
$ws.Range("E204").Value = "MIC"
$ws.Range("F204").Value = "3: 5275"
$ws.Range("G204").Value = "3: 5278"
$ws.Range("H204").Value = 0
$ws.Range("I204").Value = "_x0005_256"
$ws.Range("J204").Value = 4
$ws.Range("K204").Value = 0.0090224207154779631
$ws.Range("L204").Value = "chen"
$ws.Range("M204").Value = "1/30/19 17:59:35"

# Row 205
$ws.Range("A205").Value = "●"
$ws.Range("D205").Value = "23107"
$ws.Range("E205").Value = "MIC"
$ws.Range("F205").Value = "3: 5319"
$ws.Range("G205").Value = "3: 5322"
$ws.Range("H205").Value = 0
$ws.Range("I205").Value = "_x0005_256"
$ws.Range("J205").Value = 4
$ws.Range("K205").Value = 0.0090224207154779631
$ws.Range("L205").Value = "chen"
$ws.Range("M205").Value = "1/30/19 17:59:44"

# Row 206
$ws.Range("A206").Value = "●"
$ws.Range("D206").Value = "23107"
$ws.Range("E206").Value = "MIC"
$ws.Range("F206").Value = "3: 5356"
$ws.Range("G206").Value = "3: 5359"
$ws.Range("H206").Value = 0
$ws.Range("I206").Value = "_x0005_256"
$ws.Range("J206").Value = 4
$ws.Range("K206").Value = 0.0090224207154779631
$ws.Range("L206").Value = "chen"
$ws.Range("M206").Value = "1/30/19 17:59:50"

# Row 207
$ws.Range("A207").Value = "●"
$ws.Range("D207").Value = "23107"
$ws.Range("E207").Value = "MIC"
$ws.Range("F207").Value = "3: 5395"
$ws.Range("G207").Value = "3: 5398"
$ws.Range("H207").Value = 0
$ws.Range("I207").Value = "_x0005_256"
$ws.Range("J207").Value = 4
$ws.Range("K207").Value = 0.0090224207154779631
$ws.Range("L207").Value = "chen"
$ws.Range("M207").Value = "1/30/19 17:59:55"

# Row 208
$ws.Range("A208").Value = "●"
$ws.Range("D208").Value = "23107"
$ws.Range("E208").Value = "MIC"
$ws.Range("F208").Value = "3: 5433"
$ws.Range("G208").Value = "3: 5436"
$ws.Range("H208").Value = 0
$ws.Range("I208").Value = "_x0005_256"
$ws.Range("J208").Value = 4
$ws.Range("K208").Value = 0.0090224207154779631
$ws.Range("L208").Value = "chen"
$ws.Range("M208").Value = "1/30/19 18:00:00"

# Row 209
$ws.Range("A209").Value = "●"
$ws.Range("D209").Value = "23107"
$ws.Range("E209").Value = "MIC"
$ws.Range("F209").Value = "3: 5504"
$ws.Range("G209").Value = "3: 5506"
$ws.Range("H209").Value = 0
$ws.Range("I209").Value = "_x0005_32"
$ws.Range("J209").Value = 3
$ws.Range("K209").Value = 0.0067668155366084719
$ws.Range("L209").Value = "chen"
$ws.Range("M209").Value = "1/30/19 18:00:09"

# Row 210
$ws.Range("A210").Value = "●"
$ws.Range("D210").Value = "23107"
$ws.Range("E210").Value = "MIC"
$ws.Range("F210").Value = "3: 5530"
$ws.Range("G210").Value = "3: 5538"
$ws.Range("H210").Value = 0
$ws.Range("I210").Value = "Meropenem"
$ws.Range("J210").Value = 9
$ws.Range("K210").Value = 0.020300446609825416
$ws.Range("L210").Value = "chen"
$ws.Range("M210").Value = "1/30/19 18:00:15"

# Row 211
$ws.Range("A211").Value = "●"
$ws.Range("D211").Value = "23107"
$ws.Range("E211").Value = "MIC"
$ws.Range("F211").Value = "3: 5582"
$ws.Range("G211").Value = "3: 5584"
$ws.Range("H211").Value = 0
$ws.Range("I211").Value = "_x0005_32"
$ws.Range("J211").Value = 3
$ws.Range("K211").Value = 0.0067668155366084719
$ws.Range("L211").Value = "chen"
$ws.Range("M211").Value = "1/30/19 18:00:25"

# Row 212
$ws.Range("A212").Value = "●"
$ws.Range("D212").Value = "23107"
$ws.Range("E212").Value = "MIC"
$ws.Range("F212").Value = "3: 5239"
$ws.Range("G212").Value = "3: 5242"
$ws.Range("H212").Value = 0
$ws.Range("I212").Value = "_x0005_256"
$ws.Range("J212").Value = 4
$ws.Range("K212").Value = 0.0090224207154779631
$ws.Range("L212").Value = "chen"
$ws.Range("M212").Value = "1/30/19 18:00:54"

# Row 213
$ws.Range("A213").Value = "●"
$ws.Range("D213").Value = "23107"
$ws.Range("E213").Value = "MIC"
$ws.Range("F213").Value = "3: 5204"
$ws.Range("G213").Value = "3: 5207"
$ws.Range("H213").Value = 0
$ws.Range("I213").Value = "_x0005_256"
$ws.Range("J213").Value = 4
$ws.Range("K213").Value = 0.0090224207154779631
$ws.Range("L213").Value = "chen"
$ws.Range("M213").Value = "1/30/19 18:00:59"

# Row 214
$ws.Range("A214").Value = "●"
$ws.Range("D214").Value = "23107"
$ws.Range("E214").Value = "MIC"
$ws.Range("F214").Value = "3: 5167"
$ws.Range("G214").Value = "3: 5170"
$ws.Range("H214").Value = 0
$ws.Range("I214").Value = "_x0005_256"
$ws.Range("J214").Value = 4
$ws.Range("K214").Value = 0.0090224207154779631
$ws.Range("L214").Value = "chen"
$ws.Range("M214").Value = "1/30/19 18:01:08"
